$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "47.403.38"
$ws.Range("E2").Value = "  +3.11%  "
$ws.Range("D3").Value = "2.504.93"
$ws.Range("E3").Value = "  +2.49%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "'324.73"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.97%  "
$ws.Range("D6").Value = "'110.09"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.53%  "
$ws.Range("E7").Value = "  +1.35%  "
$ws.Range("D9").Value = "'0.538"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.49%  "
$ws.Range("D10").Value = "'39.42"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +9.99%  "
$ws.Range("E11").Value = "  +1.38%  "
$ws.Range("E12").Value = "  +0.86%  "
$ws.Range("D13").Value = "'18.55"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.84%  "
$ws.Range("D14").Value = "'7.22"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.56%  "
$ws.Range("D15").Value = "2.897.26"
$ws.Range("E15").Value = "  +2.39%  "
$ws.Range("D16").Value = "2.508.49"
$ws.Range("E16").Value = "  +3.35%  "
$ws.Range("D17").Value = "'0.861"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.57%  "
$ws.Range("D18").Value = "47.353.51"
$ws.Range("E18").Value = "  +3.39%  "
$ws.Range("D19").Value = "'12.90"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.44%  "
$ws.Range("D20").Value = "'6.68"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.29%  "
$ws.Range("D21").Value = "0.0₃0942"
$ws.Range("E21").Value = "  +1.13%  "
$ws.Range("D22").Value = "'2.71"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +14.21%  "
$ws.Range("D23").Value = "'70.61"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.84%  "
$ws.Range("D24").Value = "'248.91"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.96%  "
$ws.Range("E25").Value = "  +3.85%  "
$ws.Range("D26").Value = "'26.13"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.08%  "
$ws.Range("D27").Value = "'1.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.02%  "
$ws.Range("B28").Value = "Cosmos"
$ws.Range("C28").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D28").Value = "'10.08"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.29%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").Value = "'2.21"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.03%  "
$ws.Range("D30").Value = "'35.50"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.75%  "
$ws.Range("D31").Value = "'0.138"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +9.90%  "
$ws.Range("D32").Value = "'49.94"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.42%  "
$ws.Range("D33").Value = "'20.03"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.07%  "
$ws.Range("D34").Value = "'5.43"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.54%  "
$ws.Range("E35").Value = "  +5.33%  "
$ws.Range("E36").Value = "  +0.23%  "
$ws.Range("E37").Value = "  +6.11%  "
$ws.Range("D38").Value = "'4.71"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.19%  "
$ws.Range("E39").Value = "  +2.02%  "
$ws.Range("E40").Value = "  +1.70%  "
$ws.Range("D41").Value = "'122.03"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.28%  "
$ws.Range("E42").Value = "  -1.61%  "
$ws.Range("D43").Value = "'21.37"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.75%  "
$ws.Range("D44").Value = "'0.0301"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.84%  "
$ws.Range("D45").Value = "2.002.02"
$ws.Range("E45").Value = "  +2.30%  "
$ws.Range("D46").Value = "'3.13"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +5.97%  "
$ws.Range("E47").Value = "  -2.20%  "
$ws.Range("E48").Value = "  -3.68%  "
$ws.Range("D49").Value = "'9.06"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.26%  "
$ws.Range("D50").Value = "'5.23"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +5.12%  "
$ws.Range("B51").Value = "MultiversX"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D51").Value = "'56.76"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.42%  "
